# Update the marking scheme / total marks on the "quiz" marksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row (row 11): marks awarded per correct answer, 3 -> 5
$ws.Range("B11").Value = 5

# Total row (row 12): total marks obtained by correct answers, 36 -> 60
$ws.Range("B12").Value = 60

# Total row (row 12), Max column: "obtained/total max marks" text, 36/84 -> 60/140
$ws.Range("E12").Value = "60/140"
